$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3300368681588889,0.295472258886636,0.3631211525990467,0.2722860630512539,0.4068499506284093,7,7,7,7,7,4,4,4,4,4,2.310258077112223,2.068305812206452,2.541848068193326,1.906002441358778,2.847949654398865,100.58425,87.68037383032676,112.350909364063,77.78486367272372,128.7766997559584,0.127699338494011,0.1038513871225602,0.151622640529654,0.08604351479746133,0.1749317536996756,0.8571964331914291,0.8112199873047428,0.8979354764039813,0.7705929268849784,0.9301073959422517),
    @(0.33,0.33,0.33,0.33,0.33,6.997208238996328,6.211708480501438,7.750615725360402,5.697852097493597,8.761324523328691,4,4,4,4,4,2.309078718868789,2.049863798565474,2.557703189368933,1.880291192172887,2.891237092698468,99.94318,92.74969224531547,107.3669139784783,87.51808929487457,115.0151989085309,0.1280655974290327,0.09771137674306317,0.1597094549971312,0.07615470438817838,0.1925575142439491,0.8559038826910592,0.806321355153669,0.8998248219988605,0.7628449515139616,0.9340592150328463),
    @(0.3299546669000036,0.2954865006110795,0.3628830664704524,0.2720329806800469,0.407098622037785,6.99827874932453,6.211970424988903,7.745037165329355,5.701645464256834,8.770715711799282,4,4,4,4,4,2.309030454651602,1.957856982297311,2.64493494335285,1.725062703407883,3.108858007767917,102.08827,87.09196012881102,116.5183351204225,76.61988320651332,138.3072911618482,0.127053027103331,0.08886379706095983,0.1661520603053158,0.06100639498027485,0.2056703573250038,0.8487558790433356,0.7772784448030349,0.9078315424462871,0.7061324731648092,0.9513991009616451),
    @(0.33,0.33,0.33,0.33,0.33,7,7,7,7,7,3.99890249252565,3.54783857102002,4.431856646837604,3.165981920593907,4.927173857001854,2.31,2.31,2.31,2.31,2.31,98.52849000000001,93.79558246075284,103.1431825745692,89.93686124676002,107.6585286999046,0.1289788574428249,0.1233500633980511,0.1343574825850218,0.1186524547236349,0.1394438571624682,0.86416766455395,0.8641670458907781,0.8641680178366636,0.8641659827486737,0.8641683396216523),
    @(0.33002046411242,0.2953763363915379,0.3632716599738005,0.2718840422252982,0.4065125934875387,6.998500463445088,6.209881634710266,7.756141471337187,5.698366000827131,8.768556060092456,3.999545593737989,3.548838938696382,4.432336156891433,3.172164494650679,4.920191628817671,2.309632089802328,1.957416985369304,2.647699099948058,1.721622927877045,3.111835247514327,102.04347,86.25076848911166,116.9307385069291,75.21182886888346,138.7914045077047,0.1273154517357958,0.08857550152371189,0.1670925343653169,0.0603878923248169,0.207527434219609,0.8486599693161534,0.7767003830988098,0.9080318014720987,0.7047214527548993,0.9519399972407894)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
